$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title: "Version 1.1" -> "Version 1.4"
#    (the revision-history paragraph further down also contains the
#    literal text "Version 1.1" and must stay untouched, so locate the
#    title paragraph specifically -- the one whose entire text is just
#    "Version 1.1" -- and scope the Find/Replace to that paragraph only.)
# ---------------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Version 1.1") {
        $titlePara = $p
        break
    }
}
if ($titlePara -eq $null) {
    throw "Could not find the 'Version 1.1' title paragraph"
}

$titleRange = $titlePara.Range
$replaced = $titleRange.Find.Execute("1.1", $false, $false, $false, $false, $false, $true, 0, $false, "1.4", 2)
if (-not $replaced) {
    throw "Failed to replace '1.1' with '1.4' in the title"
}

# ---------------------------------------------------------------------
# 2. Revision history: append a new entry after the
#    "Version 1.1 ... camera bugfix" line documenting the 1.4 release.
# ---------------------------------------------------------------------
$historyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Version 1.1　…　カメラ起動中に転送が停止する不具合を修正") {
        $historyPara = $p
    }
}
if ($historyPara -eq $null) {
    throw "Could not find the revision-history paragraph to append after"
}

$historyRange = $historyPara.Range
$historyRange.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.Text = "Version 1.4　…　カメラの安定性を改善"
